# Daily attendance processing - 2026-01-16 21:36:03
# Rotate the "Recorded By" (column G) comma-separated list of names by
# moving the first entry to the end, for every row that has more than
# one recorder listed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $value = $cell.Value()

    if ($value -ne $null -and $value -like "*,*") {
        $parts = $value -split ",\s*"
        $rotated = ($parts[1..($parts.Length - 1)] + $parts[0]) -join ", "
        $cell.Value = $rotated
    }
}
